$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "39.860.04"
Set-TextValue "E2" "  -3.09%  "
Set-TextValue "D3" "2.334.33"
Set-TextValue "E3" "  -3.85%  "
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "309.28"
Set-TextValue "E5" "  -2.75%  "
Set-TextValue "D6" "83.80"
Set-TextValue "E6" "  -6.31%  "
Set-TextValue "E7" "  -2.23%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "E9" "  -3.95%  "
Set-TextValue "E10" "  -4.20%  "
Set-TextValue "D11" "29.71"
Set-TextValue "E11" "  -7.36%  "
Set-TextValue "E12" "  +0.55%  "
Set-TextValue "D13" "2.694.19"
Set-TextValue "E13" "  -3.91%  "
Set-TextValue "E14" "  -5.33%  "
Set-TextValue "D15" "14.68"
Set-TextValue "E15" "  -6.15%  "
Set-TextValue "D16" "2.353.98"
Set-TextValue "E16" "  -2.92%  "
Set-TextValue "E17" "  -3.20%  "
Set-TextValue "D18" "39.793.96"
Set-TextValue "E19" "  -3.35%  "
Set-TextValue "E20" "  -4.21%  "
Set-TextValue "D21" "67.92"
Set-TextValue "E21" "  -6.14%  "
Set-TextValue "D22" "10.48"
Set-TextValue "E22" "  -5.19%  "
Set-TextValue "D23" "234.12"
Set-TextValue "E23" "  -0.45%  "
Set-TextValue "E24" "  -6.26%  "
Set-TextValue "E25" "  -0.02%  "
Set-TextValue "E26" "  -3.38%  "
Set-TextValue "E27" "  -3.31%  "
Set-TextValue "E28" "  -1.51%  "
Set-TextValue "D29" "9.21"
Set-TextValue "E29" "  -4.38%  "
Set-TextValue "D30" "33.84"
Set-TextValue "E30" "  -2.09%  "
Set-TextValue "D31" "152.68"
Set-TextValue "E31" "  -3.60%  "
Set-TextValue "E32" "  -0.17%  "
Set-TextValue "E33" "  -3.92%  "
Set-TextValue "E34" "  -0.94%  "
Set-TextValue "E35" "  -4.45%  "
Set-TextValue "E36" "  -0.89%  "
Set-TextValue "D37" "2.76"
Set-TextValue "E37" "  -6.38%  "
Set-TextValue "D38" "0.0977"
Set-TextValue "E38" "  -3.11%  "
Set-TextValue "D39" "15.44"
Set-TextValue "E39" "  -9.48%  "
Set-TextValue "E40" "  -5.20%  "
Set-TextValue "E41" "  -3.50%  "
Set-TextValue "D42" "1.972.59"
Set-TextValue "E42" "  -1.13%  "
Set-TextValue "E43" "  -3.35%  "
Set-TextValue "E44" "  -5.03%  "
Set-TextValue "D45" "17.35"
Set-TextValue "E45" "  -6.39%  "
Set-TextValue "D46" "9.39"
Set-TextValue "E46" "  -1.49%  "
Set-TextValue "E47" "  -8.36%  "
Set-TextValue "D48" "2.558.72"
Set-TextValue "E48" "  -4.03%  "
Set-TextValue "D49" "92.15"
Set-TextValue "E49" "  -2.82%  "
Set-TextValue "D50" "69.74"
Set-TextValue "E50" "  -5.09%  "
Set-TextValue "E51" "  -4.60%  "
